$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the completed "Entity Animations" task (old row 59); rows below shift up by one.
$ws.Rows.Item(59).Delete()

# Developers updated on two tasks that now also involve Alpi.
$ws.Range("E61").Value = "Alpi, Parrinello"
$ws.Range("E62").Value = "Alpi, Foschini"

# Record Sprint 4 effort for the "Add button for start the match" task.
$ws.Range("J52").Value = 2

# Record Sprint 5 effective effort for the fifth sprint review.
$ws.Range("K53").Value = 5
$ws.Range("K54").Value = 4
$ws.Range("K55").Value = 4
$ws.Range("K56").Value = 3
$ws.Range("K59").Value = 4
$ws.Range("K60").Value = 6
$ws.Range("K61").Value = 5
$ws.Range("K62").Value = 4

# Add the new "Add fade animation" task under the "Fix bugs" backlog item.
$ws.Range("B62:F62").Copy()
$ws.Range("B63:F63").PasteSpecial(-4122)
$ws.Range("C63").Value = "Add fade animation"
$ws.Range("E63").Value = "Alpi"
$ws.Range("F63").Value = 4
$ws.Range("K63").Value = 4

# Leave the selection on the newly added cell, like the author did.
$ws.Range("K63").Select()
